$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row (row 11) with the updated results data
$ws.Range("A11").Value = 3
$ws.Range("B11").Value = 50
$ws.Range("C11").Value = 0.005889368057250976
$ws.Range("D11").Value = 1111.130075874329
$ws.Range("E11").Value = 1200.114905548096
$ws.Range("F11").Value = "-"
$ws.Range("G11").Value = 75.69570599999999
$ws.Range("H11").Value = 135.962256
$ws.Range("I11").Value = 489.28846
$ws.Range("J11").Value = "-"
